# Swap the deck's colour theme from "Integral" (green) back to the
# stock "Office Theme" (blue) palette. The presentation's single active
# DrawingML theme -- the one bound to the slide master / presentation
# (persisted as ppt/theme/theme2.xml) -- has its 12 theme colours
# rewritten, index-for-index, to the default "Office Theme" values.

function HexToBgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

$officeTheme = @(
    "000000",   # 1  dk1
    "FFFFFF",   # 2  lt1
    "44546A",   # 3  dk2
    "E7E6E6",   # 4  lt2
    "5B9BD5",   # 5  accent1
    "ED7D31",   # 6  accent2
    "A5A5A5",   # 7  accent3
    "FFC000",   # 8  accent4
    "4472C4",   # 9  accent5
    "70AD47",   # 10 accent6
    "0563C1",   # 11 hlink
    "954F72"    # 12 folHlink
)

for ($i = 1; $i -le $officeTheme.Length; $i++) {
    $colors.Item($i).RGB = HexToBgr $officeTheme[$i - 1]
}
